# The "Haba" (fava bean) price sheet gets a new weekly price observation.
# A new row is inserted right before the current row 14, which shifts all
# subsequent rows (old 14-49) down by one (to 15-50), and the new row 14 is
# populated with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14 - this shifts rows 14..49 down to 15..50
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with the new price observation
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44707
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112026
$ws.Range("G14").Value = "Haba"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13500
$ws.Range("N14").Value = "`$/saco 25 kilos"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 540
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
